$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Bonesplitter', ['{1}', 'Artifact — Equipment', 'Equipped creature gets +2/+0.', 'Equip {1}'])"
$ws.Range("A3").Value = "('Elvish Aberration', ['{5}{G}', 'Creature — Elf Mutant', '{T}: Add {G}{G}{G}.', 'Forestcycling {2} ({2}, Discard this card: Search your library for a Forest card, reveal it, put it into your hand, then shuffle your library.)', '4/5'])"
$ws.Range("A4").Value = "('Forest', ['Basic Land — Forest', '({T}: Add {G}.)'])"
$ws.Range("A5").Value = "('Island', ['Basic Land — Island', '({T}: Add {U}.)'])"
$ws.Range("A6").Value = "('Mountain', ['Basic Land — Mountain', '({T}: Add {R}.)'])"
$ws.Range("A7").Value = "('Plains', ['Basic Land — Plains', '({T}: Add {W}.)'])"
$ws.Range("A8").Value = "('Skirk Marauder', ['{1}{R}', 'Creature — Goblin', 'Morph {2}{R} (You may cast this card face down as a 2/2 creature for {3}. Turn it face up any time for its morph cost.)', 'When Skirk Marauder is turned face up, it deals 2 damage to any target.', '2/1'])"
$ws.Range("A9").Value = "('Swamp', ['Basic Land — Swamp', '({T}: Add {B}.)'])"

$ws.Range("A10:A33").ClearContents()
